$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1554.5
$ws.Range("I99").Value = 252
$ws.Range("J99").Value = 3725.3333
$ws.Range("K99").Value = 756
$ws.Range("L99").Value = 11175.9999
$ws.Range("M99").Value = 742
$ws.Range("N99").Value = -14171.9999

$ws.Range("H137").Value = 3494.8125
$ws.Range("I137").Value = 2605.76
$ws.Range("J137").Value = 6670
$ws.Range("K137").Value = 7817.280000000001
$ws.Range("L137").Value = 20010
$ws.Range("M137").Value = -5267.280000000001
$ws.Range("N137").Value = -25110

$ws.Range("H138").Value = 1825.717
$ws.Range("I138").Value = 1429.8611
$ws.Range("K138").Value = 4289.5833
$ws.Range("M138").Value = 850.4166999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 636043.8
$ws.Range("I32").Value = 808713.8
$ws.Range("J32").Value = 17309.5
$ws.Range("K32").Value = 808713.8
$ws.Range("L32").Value = 17309.5
$ws.Range("M32").Value = -808426.8
$ws.Range("N32").Value = -17883.5

$ws.Range("H61").Value = 2765.7144
$ws.Range("I61").Value = 1721.6154
$ws.Range("K61").Value = 1721.6154
$ws.Range("M61").Value = -1509.6154

$ws.Range("H74").Value = 1435.238
$ws.Range("I74").Value = 1023.1429
$ws.Range("J74").Value = 2259.4285
$ws.Range("K74").Value = 1023.1429
$ws.Range("L74").Value = 2259.4285
$ws.Range("M74").Value = -149.1429000000001
$ws.Range("N74").Value = -4007.4285

$ws.Range("H77").Value = 1435.238
$ws.Range("I77").Value = 1023.1429
$ws.Range("J77").Value = 2259.4285
$ws.Range("K77").Value = 5115.7145
$ws.Range("L77").Value = 11297.1425
$ws.Range("M77").Value = -747.7145
$ws.Range("N77").Value = -20033.1425

$ws.Range("H97").Value = 1017.375
$ws.Range("I97").Value = 1017.375
$ws.Range("K97").Value = 1017.375
$ws.Range("M97").Value = -521.375

$ws.Range("H110").Value = 1256.8889
$ws.Range("I110").Value = 1184.4
$ws.Range("J110").Value = 1347.5
$ws.Range("K110").Value = 1184.4
$ws.Range("L110").Value = 1347.5
$ws.Range("M110").Value = 860.5999999999999
$ws.Range("N110").Value = -5437.5

$ws.Range("H112").Value = 79387
$ws.Range("J112").Value = 79387
$ws.Range("L112").Value = 79387
$ws.Range("N112").Value = -82341

$ws.Range("H132").Value = 3540.2678
$ws.Range("I132").Value = 2566.9023
$ws.Range("K132").Value = 7700.706900000001
$ws.Range("M132").Value = -5170.706900000001

$ws.Range("H136").Value = 2765.7144
$ws.Range("I136").Value = 1721.6154
$ws.Range("K136").Value = 5164.8462
$ws.Range("M136").Value = -2614.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1522.1333
$ws.Range("I20").Value = 1444.931
$ws.Range("K20").Value = 1444.931
$ws.Range("M20").Value = -1197.931

$ws.Range("H134").Value = 2796.6875
$ws.Range("I134").Value = 2395.5833
$ws.Range("K134").Value = 7186.749899999999
$ws.Range("M134").Value = -4651.749899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7168.8477
$ws.Range("I31").Value = 1350.409
$ws.Range("J31").Value = 12502.417
$ws.Range("K31").Value = 1350.409
$ws.Range("L31").Value = 12502.417
$ws.Range("M31").Value = -1055.409
$ws.Range("N31").Value = -13092.417

$ws.Range("H34").Value = 7168.8477
$ws.Range("I34").Value = 1350.409
$ws.Range("J34").Value = 12502.417
$ws.Range("K34").Value = 1350.409
$ws.Range("L34").Value = 12502.417
$ws.Range("M34").Value = -1148.409
$ws.Range("N34").Value = -12906.417

$ws.Range("H58").Value = 1216
$ws.Range("I58").Value = 1209.1428
$ws.Range("J58").Value = 1264
$ws.Range("K58").Value = 1209.1428
$ws.Range("L58").Value = 1264
$ws.Range("M58").Value = -1006.1428
$ws.Range("N58").Value = -1670

$ws.Range("H132").Value = 5557280
$ws.Range("I132").Value = 1429.08
$ws.Range("J132").Value = 33336534
$ws.Range("K132").Value = 4287.24
$ws.Range("L132").Value = 100009602
$ws.Range("M132").Value = -1757.24
$ws.Range("N132").Value = -100014662

$ws.Range("H134").Value = 4422.4443
$ws.Range("I134").Value = 2542.4285
$ws.Range("J134").Value = 11002.5
$ws.Range("K134").Value = 7627.2855
$ws.Range("L134").Value = 33007.5
$ws.Range("M134").Value = -5092.2855
$ws.Range("N134").Value = -38077.5

$ws.Range("H136").Value = 1216
$ws.Range("I136").Value = 1209.1428
$ws.Range("J136").Value = 1264
$ws.Range("K136").Value = 3627.4284
$ws.Range("L136").Value = 3792
$ws.Range("M136").Value = -1077.4284
$ws.Range("N136").Value = -8892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 3199.8
$ws.Range("J21").Value = 3250
$ws.Range("L21").Value = 9750
$ws.Range("N21").Value = -10096

$ws.Range("H92").Value = 598.4
$ws.Range("I92").Value = 596
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 1788
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = -540
$ws.Range("N92").Value = -4296

$ws.Range("H113").Value = 1231.3889
$ws.Range("I113").Value = 641
$ws.Range("J113").Value = 2412.1667
$ws.Range("K113").Value = 1923
$ws.Range("L113").Value = 7236.500100000001
$ws.Range("M113").Value = 247
$ws.Range("N113").Value = -11576.5001

$ws.Range("H122").Value = 6814.5
$ws.Range("I122").Value = 469.91666
$ws.Range("J122").Value = 25848.25
$ws.Range("K122").Value = 4229.24994
$ws.Range("L122").Value = 232634.25
$ws.Range("M122").Value = -1779.24994
$ws.Range("N122").Value = -237534.25

$ws.Range("H124").Value = 2175.5833
$ws.Range("I124").Value = 500
$ws.Range("J124").Value = 2327.9092
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 6983.7276
$ws.Range("M124").Value = 3410
$ws.Range("N124").Value = -16803.7276

$ws.Range("H125").Value = 2007.7
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2007.7
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 6023.1
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -15863.1

$ws.Range("H140").Value = 1476.6
$ws.Range("I140").Value = 989.9
$ws.Range("K140").Value = 2969.7
$ws.Range("M140").Value = 2210.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5010.4
$ws.Range("I70").Value = 4835.3335
$ws.Range("J70").Value = 5360.533
$ws.Range("K70").Value = 4835.3335
$ws.Range("L70").Value = 5360.533
$ws.Range("M70").Value = -4565.3335
$ws.Range("N70").Value = -5900.533

$ws.Range("H73").Value = 5010.4
$ws.Range("I73").Value = 4835.3335
$ws.Range("J73").Value = 5360.533
$ws.Range("K73").Value = 4835.3335
$ws.Range("L73").Value = 5360.533
$ws.Range("M73").Value = -3899.3335
$ws.Range("N73").Value = -7232.533

$ws.Range("H132").Value = 2151.7632
$ws.Range("I132").Value = 1805.7407
$ws.Range("J132").Value = 3001.0908
$ws.Range("K132").Value = 5417.2221
$ws.Range("L132").Value = 9003.2724
$ws.Range("M132").Value = -2887.2221
$ws.Range("N132").Value = -14063.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7667.1113
$ws.Range("I61").Value = 7250.6665
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 7250.6665
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -7048.6665
$ws.Range("N61").Value = -8904

$ws.Range("H113").Value = 7667.1113
$ws.Range("I113").Value = 7250.6665
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 7250.6665
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -5080.6665
$ws.Range("N113").Value = -12840

$ws.Range("H132").Value = 3049.375
$ws.Range("I132").Value = 2121.6
$ws.Range("J132").Value = 4595.6665
$ws.Range("K132").Value = 6364.799999999999
$ws.Range("L132").Value = 13786.9995
$ws.Range("M132").Value = -3834.799999999999
$ws.Range("N132").Value = -18846.9995

$ws.Range("H136").Value = 9806452
$ws.Range("I136").Value = 2153.923
$ws.Range("K136").Value = 6461.768999999999
$ws.Range("M136").Value = -3911.768999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11907409
$ws.Range("I132").Value = 2338.6667
$ws.Range("J132").Value = 33336536
$ws.Range("K132").Value = 7016.000100000001
$ws.Range("L132").Value = 100009608
$ws.Range("M132").Value = -4486.000100000001
$ws.Range("N132").Value = -100014668

$ws.Range("H136").Value = 2457.5715
$ws.Range("I136").Value = 2198.0344
$ws.Range("K136").Value = 6594.1032
$ws.Range("M136").Value = -4044.1032
